# Scheduled runner update: refresh Leve profit calculations across sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) with latest
# market data for ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 397
$ws.Range("J41").Value = 658.8
$ws.Range("L41").Value = 658.8
$ws.Range("N41").Value = -1538.8
# Row 107
$ws.Range("H107").Value = 629.13635
$ws.Range("J107").Value = 496.33334
$ws.Range("L107").Value = 496.33334
$ws.Range("N107").Value = -4336.33334
# Row 137
$ws.Range("H137").Value = 988.8
$ws.Range("I137").Value = 773.2353000000001
$ws.Range("K137").Value = 2319.7059
$ws.Range("M137").Value = 230.2941000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 623.8461
$ws.Range("I2").Value = 634.1667
$ws.Range("K2").Value = 634.1667
$ws.Range("M2").Value = -521.1667
# Row 45
$ws.Range("H45").Value = 1964.875
$ws.Range("I45").Value = 1997.4286
$ws.Range("K45").Value = 1997.4286
$ws.Range("M45").Value = -1620.4286
# Row 61
$ws.Range("H61").Value = 3346.5
$ws.Range("I61").Value = 3346.5
$ws.Range("K61").Value = 3346.5
$ws.Range("M61").Value = -3134.5
# Row 74
$ws.Range("H74").Value = 1494.8889
$ws.Range("I74").Value = 1494.8889
$ws.Range("K74").Value = 1494.8889
$ws.Range("M74").Value = -620.8888999999999
# Row 77
$ws.Range("H77").Value = 1494.8889
$ws.Range("I77").Value = 1494.8889
$ws.Range("K77").Value = 7474.4445
$ws.Range("M77").Value = -3106.4445
# Row 97
$ws.Range("H97").Value = 511.85715
$ws.Range("I97").Value = 480.5
$ws.Range("J97").Value = 700
$ws.Range("K97").Value = 480.5
$ws.Range("L97").Value = 700
$ws.Range("M97").Value = 15.5
$ws.Range("N97").Value = -1692
# Row 116
$ws.Range("H116").Value = 623.8461
$ws.Range("I116").Value = 634.1667
$ws.Range("K116").Value = 634.1667
$ws.Range("M116").Value = 1659.8333
# Row 122
$ws.Range("H122").Value = 1849.75
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1772.2142
$ws.Range("I132").Value = 1567.5834
$ws.Range("K132").Value = 4702.7502
$ws.Range("M132").Value = -2172.7502
# Row 136
$ws.Range("H136").Value = 3346.5
$ws.Range("I136").Value = 3346.5
$ws.Range("K136").Value = 10039.5
$ws.Range("M136").Value = -7489.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 623.8461
$ws.Range("I3").Value = 634.1667
$ws.Range("K3").Value = 634.1667
$ws.Range("M3").Value = -520.1667
# Row 20
$ws.Range("H20").Value = 1082
$ws.Range("J20").Value = 1465.6666
$ws.Range("L20").Value = 1465.6666
$ws.Range("N20").Value = -1959.6666
# Row 94
$ws.Range("H94").Value = 1271.1177
$ws.Range("I94").Value = 940.93335
$ws.Range("J94").Value = 3747.5
$ws.Range("K94").Value = 940.93335
$ws.Range("L94").Value = 3747.5
$ws.Range("M94").Value = -489.93335
$ws.Range("N94").Value = -4649.5
# Row 99
$ws.Range("H99").Value = 3135.4644
$ws.Range("I99").Value = 3252.0833
$ws.Range("J99").Value = 2435.75
$ws.Range("K99").Value = 3252.0833
$ws.Range("L99").Value = 2435.75
$ws.Range("M99").Value = -1754.0833
$ws.Range("N99").Value = -5431.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2643.875
$ws.Range("I31").Value = 2004.2667
$ws.Range("J31").Value = 3709.889
$ws.Range("K31").Value = 2004.2667
$ws.Range("L31").Value = 3709.889
$ws.Range("M31").Value = -1709.2667
$ws.Range("N31").Value = -4299.889
# Row 34
$ws.Range("H34").Value = 2643.875
$ws.Range("I34").Value = 2004.2667
$ws.Range("J34").Value = 3709.889
$ws.Range("K34").Value = 2004.2667
$ws.Range("L34").Value = 3709.889
$ws.Range("M34").Value = -1802.2667
$ws.Range("N34").Value = -4113.889
# Row 94
$ws.Range("H94").Value = 1618.4
$ws.Range("J94").Value = 1648.5
$ws.Range("L94").Value = 1648.5
$ws.Range("N94").Value = -2550.5
# Row 96
$ws.Range("H96").Value = 20984.8
$ws.Range("J96").Value = 20984.8
$ws.Range("L96").Value = 20984.8
$ws.Range("N96").Value = -26476.8
# Row 141
$ws.Range("H141").Value = 519999
$ws.Range("J141").Value = 657141.7
$ws.Range("L141").Value = 657141.7
$ws.Range("N141").Value = -667501.7

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 3996.75
$ws.Range("J131").Value = 4995.6665
$ws.Range("L131").Value = 14986.9995
$ws.Range("N131").Value = -25066.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1548
$ws.Range("J97").Value = 3900
$ws.Range("L97").Value = 3900
$ws.Range("N97").Value = -4892
# Row 103
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
# Row 111
$ws.Range("H111").Value = 75000
$ws.Range("J111").Value = 75000
$ws.Range("L111").Value = 75000
$ws.Range("N111").Value = -81134

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4185.846
$ws.Range("I61").Value = 4268.5557
$ws.Range("J61").Value = 3999.75
$ws.Range("K61").Value = 4268.5557
$ws.Range("L61").Value = 3999.75
$ws.Range("M61").Value = -4066.5557
$ws.Range("N61").Value = -4403.75
# Row 93
$ws.Range("H93").Value = 6000
$ws.Range("I93").Value = 6000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 6000
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("M93").Value = -4752
# Row 100
$ws.Range("H100").Value = 2850
$ws.Range("I100").Value = 2825.1667
$ws.Range("K100").Value = 2825.1667
$ws.Range("M100").Value = -2284.1667
# Row 113
$ws.Range("H113").Value = 4185.846
$ws.Range("I113").Value = 4268.5557
$ws.Range("J113").Value = 3999.75
$ws.Range("K113").Value = 4268.5557
$ws.Range("L113").Value = 3999.75
$ws.Range("M113").Value = -2098.5557
$ws.Range("N113").Value = -8339.75
# Row 141
$ws.Range("H141").Value = 100000
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 557.2
$ws.Range("I107").Value = 312.5
$ws.Range("J107").Value = 924.25
$ws.Range("K107").Value = 937.5
$ws.Range("L107").Value = 2772.75
$ws.Range("M107").Value = 982.5
$ws.Range("N107").Value = -6612.75
# Row 113
$ws.Range("H113").Value = 369.4
$ws.Range("I113").Value = 300.33334
$ws.Range("J113").Value = 473
$ws.Range("K113").Value = 901.0000200000001
$ws.Range("L113").Value = 1419
$ws.Range("M113").Value = 1268.99998
$ws.Range("N113").Value = -5759
# Row 122
$ws.Range("H122").Value = 1134.5625
$ws.Range("I122").Value = 938
$ws.Range("K122").Value = 2814
$ws.Range("M122").Value = -364
# Row 136
$ws.Range("H136").Value = 2847.8
$ws.Range("I136").Value = 3184.75
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 9554.25
$ws.Range("L136").Value = 1500
$ws.Range("M136").Value = -7004.25
$ws.Range("N136").Value = -9600
